$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "UA" language column (D) -----------------------------------------
# Header
$ws.Range("D2").Value = "UA"

# Body rows, written in the same order the original author typed them in
# (this reproduces the exact shared-string append order of the source file).
$ws.Range("D3").Value  = "Контакти"
$ws.Range("D4").Value  = "Фільтрувати контакти"
$ws.Range("D5").Value  = "Перетворити в xml"
$ws.Range("D6").Value  = "Ім'я"
$ws.Range("D7").Value  = "Телефон"
$ws.Range("D8").Value  = "Група"
$ws.Range("D9").Value  = "Редагувати контакт"
$ws.Range("D10").Value = "Видалити контакт"
$ws.Range("D11").Value = "Додати контакт"
$ws.Range("D13").Value = "Темна тема"
$ws.Range("D14").Value = "Мова"
$ws.Range("D12").Value = "Світла тема"

# Match column D's formatting (header style + bordered body style, including
# the still-empty rows 15-21) to column C's.
$ws.Range("C2:C21").Copy()
$ws.Range("D2:D21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths: B,C,D all become one uniform width ---------------------
$ws.Range("B2:D21").ColumnWidth = 24.45

# --- Selection moves to C14 (matches the saved cursor position) -----------
$ws.Range("C14").Select()
